$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Add the new "Mutación 0,00" series in column U, mirroring the
# existing R/S/T ("Mutación 0,01" / "0,05" / "0,09") helper columns
# that already feed the three box-whisker charts.
# -----------------------------------------------------------------

# Block 1 (rows 6-11) - feeds first chart
$ws.Range("U6").Value = "Mutación 0,00"
$ws.Range("U7").Value = 20526
$ws.Range("U8").Value = 20349
$ws.Range("U9").Value = 20576
$ws.Range("U10").Value = 20506
$ws.Range("U11").Value = 20445

# Block 2 (rows 13-18) - feeds second chart
$ws.Range("U13").Value = "Mutación 0,00"
$ws.Range("U14").Value = 35583
$ws.Range("U15").Value = 35629
$ws.Range("U16").Value = 35622
$ws.Range("U17").Value = 35505
$ws.Range("U18").Value = 35553

# Block 3 (rows 20-25) - feeds third chart
$ws.Range("U20").Value = "Mutación 0,00"
$ws.Range("U21").Value = 4508
$ws.Range("U22").Value = 4462
$ws.Range("U23").Value = 4454
$ws.Range("U24").Value = 4498
$ws.Range("U25").Value = 4488

# -----------------------------------------------------------------
# Drop the old ad-hoc yellow fill / thick borders that used to mark
# up the E:G helper block - the sheet now relies on plain formatting
# there (matches the cleaned up styles.xml / row heights).
# -----------------------------------------------------------------
$ws.Range("E7:G11").Style = "Normal"
$ws.Range("R14:R18").Style = "Normal"
$ws.Range("R21:R25").Style = "Normal"

# Clear the explicit row heights / bottom border flags that came
# along with the old look of rows 7-11, 14-18 and 21-25.
$ws.Rows("7:11").RowHeight = $ws.Rows("6:6").RowHeight
$ws.Rows("14:18").RowHeight = $ws.Rows("6:6").RowHeight
$ws.Rows("21:25").RowHeight = $ws.Rows("6:6").RowHeight

# -----------------------------------------------------------------
# Update the view: scroll the sheet and move the active selection.
# -----------------------------------------------------------------
$ws.Range("AC62").Select()

Write-Host "done"
